# Auto-generated edit script applying numeric corrections to LeveProfit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 3000
$ws.Range("I21").Value = 3000
$ws.Range("K21").Value = 3000
$ws.Range("M21").Value = -2532
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 3000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = -2766
$ws.Range("H28").Value = 1703.1052
$ws.Range("I28").Value = 1242.2222
$ws.Range("J28").Value = 9999
$ws.Range("K28").Value = 1242.2222
$ws.Range("L28").Value = 9999
$ws.Range("M28").Value = -757.2221999999999
$ws.Range("N28").Value = -10969
$ws.Range("H40").Value = 3708341
$ws.Range("J40").Value = 6672226.5
$ws.Range("L40").Value = 6672226.5
$ws.Range("N40").Value = -6672576.5
$ws.Range("H51").Value = 14203.556
$ws.Range("I51").Value = 29966.666
$ws.Range("J51").Value = 6322
$ws.Range("K51").Value = 29966.666
$ws.Range("L51").Value = 6322
$ws.Range("M51").Value = -29482.666
$ws.Range("N51").Value = -7290
$ws.Range("H53").Value = 6945.3
$ws.Range("I53").Value = 5407.4287
$ws.Range("K53").Value = 5407.4287
$ws.Range("M53").Value = -4770.4287
$ws.Range("H64").Value = 111118264
$ws.Range("J64").Value = 8133.1665
$ws.Range("L64").Value = 8133.1665
$ws.Range("N64").Value = -8629.1665
$ws.Range("H67").Value = 111118264
$ws.Range("J67").Value = 8133.1665
$ws.Range("L67").Value = 8133.1665
$ws.Range("N67").Value = -9849.1665
$ws.Range("H68").Value = 71995
$ws.Range("J68").Value = 71995
$ws.Range("L68").Value = 71995
$ws.Range("N68").Value = -73493
$ws.Range("H71").Value = 71995
$ws.Range("J71").Value = 71995
$ws.Range("L71").Value = 215985
$ws.Range("N71").Value = -223473
$ws.Range("H74").Value = 166672000
$ws.Range("I74").Value = 300004200
$ws.Range("J74").Value = 6750
$ws.Range("K74").Value = 300004200
$ws.Range("L74").Value = 6750
$ws.Range("M74").Value = -300003264
$ws.Range("N74").Value = -8622
$ws.Range("H77").Value = 166672000
$ws.Range("I77").Value = 300004200
$ws.Range("J77").Value = 6750
$ws.Range("K77").Value = 1500021000
$ws.Range("L77").Value = 33750
$ws.Range("M77").Value = -1500016320
$ws.Range("N77").Value = -43110
$ws.Range("H112").Value = 4352.25
$ws.Range("J112").Value = 4352.25
$ws.Range("L112").Value = 13056.75
$ws.Range("N112").Value = -15272.75
$ws.Range("H125").Value = 55555960
$ws.Range("I125").Value = 71428960
$ws.Range("J125").Value = 449
$ws.Range("K125").Value = 642860640
$ws.Range("L125").Value = 4041
$ws.Range("M125").Value = -642858180
$ws.Range("N125").Value = -8961
$ws.Range("H132").Value = 1049.1666
$ws.Range("I132").Value = 1021.73773
$ws.Range("J132").Value = 1383.8
$ws.Range("K132").Value = 3065.21319
$ws.Range("L132").Value = 4151.4
$ws.Range("M132").Value = -535.2131900000004
$ws.Range("N132").Value = -9211.4
$ws.Range("H137").Value = 6785.927
$ws.Range("I137").Value = 4921.5
$ws.Range("K137").Value = 14764.5
$ws.Range("M137").Value = -12214.5
$ws.Range("H138").Value = 3172.79
$ws.Range("I138").Value = 2806.6924
$ws.Range("J138").Value = 3227.4941
$ws.Range("K138").Value = 8420.0772
$ws.Range("L138").Value = 9682.4823
$ws.Range("M138").Value = -3280.0772
$ws.Range("N138").Value = -19962.4823
$ws.Range("H141").Value = 2765.5833
$ws.Range("I141").Value = 2518.7
$ws.Range("K141").Value = 7556.099999999999
$ws.Range("M141").Value = -2376.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2647.5
$ws.Range("J13").Value = 3735.4
$ws.Range("L13").Value = 3735.4
$ws.Range("N13").Value = -4023.4
$ws.Range("H32").Value = 1813.06
$ws.Range("I32").Value = 1813.06
$ws.Range("K32").Value = 1813.06
$ws.Range("M32").Value = -1526.06
$ws.Range("H61").Value = 22227724
$ws.Range("I61").Value = 2397.2307
$ws.Range("J61").Value = 52641330
$ws.Range("K61").Value = 2397.2307
$ws.Range("L61").Value = 52641330
$ws.Range("M61").Value = -2185.2307
$ws.Range("N61").Value = -52641754
$ws.Range("H63").Value = 2049.1667
$ws.Range("J63").Value = 2999
$ws.Range("L63").Value = 2999
$ws.Range("N63").Value = -4371
$ws.Range("H66").Value = 2049.1667
$ws.Range("J66").Value = 2999
$ws.Range("L66").Value = 14995
$ws.Range("N66").Value = -21859
$ws.Range("H74").Value = 22194.49
$ws.Range("I74").Value = 29885.314
$ws.Range("K74").Value = 29885.314
$ws.Range("M74").Value = -29011.314
$ws.Range("H77").Value = 22194.49
$ws.Range("I77").Value = 29885.314
$ws.Range("K77").Value = 149426.57
$ws.Range("M77").Value = -145058.57
$ws.Range("H88").Value = 100000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 100000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 100000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -100812
$ws.Range("H91").Value = 100000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 100000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 100000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -102808
$ws.Range("H110").Value = 30313756
$ws.Range("I110").Value = 14511.125
$ws.Range("K110").Value = 14511.125
$ws.Range("M110").Value = -12466.125
$ws.Range("H115").Value = 40000
$ws.Range("I115").Value = 40000
$ws.Range("K115").Value = 40000
$ws.Range("M115").Value = -38433
$ws.Range("H132").Value = 5598.125
$ws.Range("I132").Value = 2987.923
$ws.Range("K132").Value = 8963.769
$ws.Range("M132").Value = -6433.769
$ws.Range("H136").Value = 22227724
$ws.Range("I136").Value = 2397.2307
$ws.Range("J136").Value = 52641330
$ws.Range("K136").Value = 7191.6921
$ws.Range("L136").Value = 157923990
$ws.Range("M136").Value = -4641.6921
$ws.Range("N136").Value = -157929090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1674.8
$ws.Range("I36").Value = 1674.8
$ws.Range("K36").Value = 1674.8
$ws.Range("M36").Value = -1140.8
$ws.Range("H86").Value = 119279.11
$ws.Range("I86").Value = 206999.8
$ws.Range("J86").Value = 9628.25
$ws.Range("K86").Value = 206999.8
$ws.Range("L86").Value = 9628.25
$ws.Range("M86").Value = -205876.8
$ws.Range("N86").Value = -11874.25
$ws.Range("H89").Value = 119279.11
$ws.Range("I89").Value = 206999.8
$ws.Range("J89").Value = 9628.25
$ws.Range("K89").Value = 1034999
$ws.Range("L89").Value = 48141.25
$ws.Range("M89").Value = -1029383
$ws.Range("N89").Value = -59373.25
$ws.Range("H105").Value = 2674.9736
$ws.Range("I105").Value = 1445.3914
$ws.Range("K105").Value = 1445.3914
$ws.Range("M105").Value = 301.6086
$ws.Range("H107").Value = 70380340
$ws.Range("I107").Value = 70380340
$ws.Range("K107").Value = 70380340
$ws.Range("M107").Value = -70378420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1199.5555
$ws.Range("I22").Value = 1212
$ws.Range("K22").Value = 1212
$ws.Range("M22").Value = -862
$ws.Range("H31").Value = 8225.714
$ws.Range("I31").Value = 3026
$ws.Range("K31").Value = 3026
$ws.Range("M31").Value = -2731
$ws.Range("H34").Value = 8225.714
$ws.Range("I34").Value = 3026
$ws.Range("K34").Value = 3026
$ws.Range("M34").Value = -2824
$ws.Range("H62").Value = 4616.143
$ws.Range("I62").Value = 2475.5
$ws.Range("J62").Value = 7470.3335
$ws.Range("K62").Value = 2475.5
$ws.Range("L62").Value = 7470.3335
$ws.Range("M62").Value = -1851.5
$ws.Range("N62").Value = -8718.3335
$ws.Range("H65").Value = 4616.143
$ws.Range("I65").Value = 2475.5
$ws.Range("J65").Value = 7470.3335
$ws.Range("K65").Value = 12377.5
$ws.Range("L65").Value = 37351.6675
$ws.Range("M65").Value = -9257.5
$ws.Range("N65").Value = -43591.6675
$ws.Range("H94").Value = 1195.3077
$ws.Range("I94").Value = 1200.7693
$ws.Range("J94").Value = 1189.8462
$ws.Range("K94").Value = 1200.7693
$ws.Range("L94").Value = 1189.8462
$ws.Range("M94").Value = -749.7692999999999
$ws.Range("N94").Value = -2091.8462
$ws.Range("H105").Value = 2977607.8
$ws.Range("I105").Value = 4202360.5
$ws.Range("K105").Value = 4202360.5
$ws.Range("M105").Value = -4200613.5
$ws.Range("H107").Value = 2078.8857
$ws.Range("I107").Value = 1643.1
$ws.Range("J107").Value = 2659.9333
$ws.Range("K107").Value = 1643.1
$ws.Range("L107").Value = 2659.9333
$ws.Range("M107").Value = 276.9000000000001
$ws.Range("N107").Value = -6499.933300000001
$ws.Range("H109").Value = 40137.5
$ws.Range("J109").Value = 40137.5
$ws.Range("L109").Value = 40137.5
$ws.Range("N109").Value = -42217.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29741312
$ws.Range("I4").Value = 31297790
$ws.Range("K4").Value = 93893370
$ws.Range("M4").Value = -93893258
$ws.Range("H37").Value = 79498.5
$ws.Range("J37").Value = 79498.5
$ws.Range("L37").Value = 238495.5
$ws.Range("N37").Value = -238719.5
$ws.Range("H80").Value = 21743492
$ws.Range("I80").Value = 29415706
$ws.Range("K80").Value = 88247118
$ws.Range("M80").Value = -88246182
$ws.Range("H83").Value = 21743492
$ws.Range("I83").Value = 29415706
$ws.Range("K83").Value = 264741354
$ws.Range("M83").Value = -264736674
$ws.Range("H92").Value = 15387434
$ws.Range("I92").Value = 3000
$ws.Range("J92").Value = 19233544
$ws.Range("K92").Value = 9000
$ws.Range("L92").Value = 57700632
$ws.Range("M92").Value = -7752
$ws.Range("N92").Value = -57703128
$ws.Range("H122").Value = 1824718.2
$ws.Range("J122").Value = 626421.8
$ws.Range("L122").Value = 5637796.2
$ws.Range("N122").Value = -5642696.2
$ws.Range("H129").Value = 333987.84
$ws.Range("I129").Value = 899.5
$ws.Range("J129").Value = 500532
$ws.Range("K129").Value = 2698.5
$ws.Range("L129").Value = 1501596
$ws.Range("M129").Value = 2301.5
$ws.Range("N129").Value = -1511596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10294.105
$ws.Range("I70").Value = 9088.9
$ws.Range("J70").Value = 11633.223
$ws.Range("K70").Value = 9088.9
$ws.Range("L70").Value = 11633.223
$ws.Range("M70").Value = -8818.9
$ws.Range("N70").Value = -12173.223
$ws.Range("H73").Value = 10294.105
$ws.Range("I73").Value = 9088.9
$ws.Range("J73").Value = 11633.223
$ws.Range("K73").Value = 9088.9
$ws.Range("L73").Value = 11633.223
$ws.Range("M73").Value = -8152.9
$ws.Range("N73").Value = -13505.223
$ws.Range("H80").Value = 3345.8
$ws.Range("I80").Value = 3068.2856
$ws.Range("K80").Value = 3068.2856
$ws.Range("M80").Value = -2070.2856
$ws.Range("H83").Value = 3345.8
$ws.Range("I83").Value = 3068.2856
$ws.Range("K83").Value = 15341.428
$ws.Range("M83").Value = -10349.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7802.727
$ws.Range("I40").Value = 6066
$ws.Range("K40").Value = 6066
$ws.Range("M40").Value = -5930
$ws.Range("H55").Value = 484.3158
$ws.Range("I55").Value = 247.22223
$ws.Range("K55").Value = 247.22223
$ws.Range("M55").Value = -74.22223
$ws.Range("H61").Value = 3281.0322
$ws.Range("I61").Value = 1873.95
$ws.Range("J61").Value = 5839.364
$ws.Range("K61").Value = 1873.95
$ws.Range("L61").Value = 5839.364
$ws.Range("M61").Value = -1671.95
$ws.Range("N61").Value = -6243.364
$ws.Range("H82").Value = 33098.25
$ws.Range("I82").Value = 64258.062
$ws.Range("J82").Value = 1938.4375
$ws.Range("K82").Value = 64258.062
$ws.Range("L82").Value = 1938.4375
$ws.Range("M82").Value = -63897.062
$ws.Range("N82").Value = -2660.4375
$ws.Range("H85").Value = 33098.25
$ws.Range("I85").Value = 64258.062
$ws.Range("J85").Value = 1938.4375
$ws.Range("K85").Value = 64258.062
$ws.Range("L85").Value = 1938.4375
$ws.Range("M85").Value = -63010.062
$ws.Range("N85").Value = -4434.4375
$ws.Range("H93").Value = 2196.7334
$ws.Range("I93").Value = 2940.2222
$ws.Range("J93").Value = 1081.5
$ws.Range("K93").Value = 2940.2222
$ws.Range("L93").Value = 1081.5
$ws.Range("M93").Value = -1692.2222
$ws.Range("N93").Value = -3577.5
$ws.Range("H100").Value = 5813.273
$ws.Range("I100").Value = 4818.385
$ws.Range("J100").Value = 7250.3335
$ws.Range("K100").Value = 4818.385
$ws.Range("L100").Value = 7250.3335
$ws.Range("M100").Value = -4277.385
$ws.Range("N100").Value = -8332.3335
$ws.Range("H113").Value = 3281.0322
$ws.Range("I113").Value = 1873.95
$ws.Range("J113").Value = 5839.364
$ws.Range("K113").Value = 1873.95
$ws.Range("L113").Value = 5839.364
$ws.Range("M113").Value = 296.05
$ws.Range("N113").Value = -10179.364
$ws.Range("H122").Value = 4107.423
$ws.Range("I122").Value = 3506.5676
$ws.Range("K122").Value = 10519.7028
$ws.Range("M122").Value = -8069.702799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1364821.9
$ws.Range("I81").Value = 1817966.4
$ws.Range("J81").Value = 5388.5
$ws.Range("K81").Value = 3635932.8
$ws.Range("L81").Value = 10777
$ws.Range("M81").Value = -3634871.8
$ws.Range("N81").Value = -12899
$ws.Range("H84").Value = 1364821.9
$ws.Range("I84").Value = 1817966.4
$ws.Range("J84").Value = 5388.5
$ws.Range("K84").Value = 18179664
$ws.Range("L84").Value = 53885
$ws.Range("M84").Value = -18174360
$ws.Range("N84").Value = -64493
$ws.Range("H107").Value = 11905835
$ws.Range("I107").Value = 700.61536
$ws.Range("K107").Value = 2101.84608
$ws.Range("M107").Value = -181.8460800000003
$ws.Range("H126").Value = 1938
$ws.Range("I126").Value = 1775
$ws.Range("J126").Value = 1978.75
$ws.Range("K126").Value = 5325
$ws.Range("L126").Value = 5936.25
$ws.Range("M126").Value = -2855
$ws.Range("N126").Value = -10876.25
$ws.Range("H132").Value = 3592.9385
$ws.Range("I132").Value = 3926.366
$ws.Range("K132").Value = 11779.098
$ws.Range("M132").Value = -9249.098
$ws.Range("H136").Value = 10981384
$ws.Range("I136").Value = 14287550
$ws.Range("J136").Value = 461764.28
$ws.Range("K136").Value = 42862650
$ws.Range("L136").Value = 1385292.84
$ws.Range("M136").Value = -42860100
$ws.Range("N136").Value = -1390392.84

Write-Host "Applied 375 cell updates and 2 clears."
